# "Pio's ERA operativo Abril-Diciembre 2025"
#
# Fernandez (row 2) left the roster; every remaining doctor shifts up one
# row, and the whole sheet's Antiguedad / Carga Teorica figures are
# recalculated for the new Abril-Diciembre 2025 operating period. The new
# recruit that lands on the now-last row (16, "Recluta1") also gets an
# updated "Fecha Ingreso a SUCA" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Fernandez's row entirely; everybody below shifts up by one.
$ws.Rows(2).Delete()

# Recalculated "Antiguedad" (C) and "Carga Teorica" (F) for every doctor
# row now that the roster/period changed.
$ws.Cells.Item(2, 3).Value = 12.09315068493151
$ws.Cells.Item(2, 6).Value = 6.215593413784655

$ws.Cells.Item(3, 3).Value = 12.09315068493151
$ws.Cells.Item(3, 6).Value = 6.215593413784655

$ws.Cells.Item(4, 3).Value = 11.58904109589041
$ws.Cells.Item(4, 6).Value = 6.272659886921544

$ws.Cells.Item(5, 3).Value = 11.58904109589041
$ws.Cells.Item(5, 6).Value = 6.272659886921544

$ws.Cells.Item(6, 3).Value = 11.00821917808219
$ws.Cells.Item(6, 6).Value = 6.338410388579262

$ws.Cells.Item(7, 3).Value = 11.00821917808219
$ws.Cells.Item(7, 6).Value = 6.338410388579262

$ws.Cells.Item(8, 3).Value = 11.00821917808219
$ws.Cells.Item(8, 6).Value = 6.338410388579262

$ws.Cells.Item(9, 3).Value = 8.586301369863014
$ws.Cells.Item(9, 6).Value = 6.61257757473692

$ws.Cells.Item(10, 3).Value = 7.753424657534246
$ws.Cells.Item(10, 6).Value = 6.706861312963085

$ws.Cells.Item(11, 3).Value = 7.553424657534246
$ws.Cells.Item(11, 6).Value = 6.729501815892393

$ws.Cells.Item(12, 3).Value = 5.753424657534246
$ws.Cells.Item(12, 6).Value = 6.933266342256174

$ws.Cells.Item(13, 3).Value = 3.917808219178082
$ws.Cells.Item(13, 6).Value = 7.141062739004625

$ws.Cells.Item(14, 3).Value = 3.16986301369863
$ws.Cells.Item(14, 6).Value = 7.225732017082725

$ws.Cells.Item(15, 3).Value = 3.16986301369863
$ws.Cells.Item(15, 6).Value = 7.225732017082725

# Row 16 is now "Recluta1" (shifted up from 17); refresh his ingreso date
# plus his Antiguedad / Carga Teorica.
$ws.Cells.Item(16, 2).Value = 45383
$ws.Cells.Item(16, 3).Value = 1.334246575342466
$ws.Cells.Item(16, 6).Value = 7.433528413831178
